$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.485269784927368
$ws.Range("B1").Value = 1.525351524353027
$ws.Range("C1").Value = 3.683701753616333
$ws.Range("D1").Value = 2.319546937942505
$ws.Range("E1").Value = 0.8427765965461731
